$d = $word.ActiveDocument

# 1. Update the date in the letter header: "July 9, 2021" -> "July 29, 2021"
$d.Content.Find.Execute("July 9, 2021", $true, $false, $false, $false, $false, `
    $true, 1, $false, "July 29, 2021", 2) | Out-Null

# 2. Split the signature block paragraph so that
#    ", Vanderbilt Memory & Alzheimer's Center" becomes its own paragraph
#    reading "Vanderbilt Memory & Alzheimer's Center" (dropping the leading
#    comma and space), keeping the same paragraph formatting.
$apostrophe = [char]0x2019
$target = ", Vanderbilt Memory & Alzheimer" + $apostrophe + "s Center"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

if ($found) {
    # Range covering just the leading ", " that needs to be removed
    $prefix = $d.Range($rng.Start, $rng.Start + 2)

    # Move the range start past the ", " so it only spans the remaining text
    $rng.Start = $rng.Start + 2

    # Break the paragraph right before "Vanderbilt Memory & Alzheimer's Center"
    $rng.InsertParagraphBefore()

    # Remove the now-orphaned ", " left behind at the end of the previous paragraph
    $prefix.Text = ""
}
